$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Academic History")
$ws.Range("F2").Value = "Monograph"
Write-Output "done"
